$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '60.598.93'
Set-TextValue 'E2' '  -1.53%  '
Set-TextValue 'D3' '2.901.75'
Set-TextValue 'E3' '  -2.42%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '526.67'
Set-TextValue 'E5' '  -2.22%  '
Set-TextValue 'D6' '142.96'
Set-TextValue 'E6' '  -4.71%  '
Set-TextValue 'E8' '  -3.68%  '
Set-TextValue 'D9' '2.909.38'
Set-TextValue 'E10' '  -5.28%  '
Set-TextValue 'D11' '5.98'
Set-TextValue 'E11' '  -2.56%  '
Set-TextValue 'D12' '0.359'
Set-TextValue 'E12' '  -2.54%  '
Set-TextValue 'D13' '3.411.78'
Set-TextValue 'E13' '  -2.45%  '
Set-TextValue 'E14' '  +2.53%  '
Set-TextValue 'D15' '60.550.91'
Set-TextValue 'E15' '  -1.73%  '
Set-TextValue 'E16' '  -4.58%  '
Set-TextValue 'D17' '2.909.18'
Set-TextValue 'E17' '  -2.32%  '
Set-TextValue 'D18' '0.0000141'
Set-TextValue 'E18' '  -3.93%  '
Set-TextValue 'E19' '  -3.71%  '
Set-TextValue 'D20' '11.59'
Set-TextValue 'E20' '  -3.72%  '
Set-TextValue 'D21' '350.82'
Set-TextValue 'E21' '  -7.47%  '
Set-TextValue 'D22' '6.52'
Set-TextValue 'E22' '  -2.50%  '
Set-TextValue 'E23' '  +0.00%  '
Set-TextValue 'E24' '  +0.92%  '
Set-TextValue 'D25' '64.60'
Set-TextValue 'E25' '  -1.65%  '
Set-TextValue 'E26' '  -3.96%  '
Set-TextValue 'E27' '  -5.56%  '
Set-TextValue 'D28' '0.997'
Set-TextValue 'E28' '  +0.10%  '
Set-TextValue 'D29' '7.81'
Set-TextValue 'E29' '  -4.83%  '
Set-TextValue 'D30' '0.0₃0848'
Set-TextValue 'E30' '  -9.49%  '
Set-TextValue 'E31' '  -0.05%  '
Set-TextValue 'E32' '  -2.69%  '
Set-TextValue 'D33' '19.58'
Set-TextValue 'E33' '  -4.12%  '
Set-TextValue 'D34' '150.89'
Set-TextValue 'E34' '  -5.23%  '
Set-TextValue 'E35' '  -6.97%  '
Set-TextValue 'E36' '  -5.97%  '
Set-TextValue 'D37' '0.996'
Set-TextValue 'E37' '  -6.84%  '
Set-TextValue 'E38' '  -5.85%  '
Set-TextValue 'D39' '37.69'
Set-TextValue 'E39' '  +0.35%  '
Set-TextValue 'E41' '  -5.34%  '
Set-TextValue 'E42' '  -5.06%  '
Set-TextValue 'E43' '  -3.39%  '
Set-TextValue 'E44' '  -1.87%  '
Set-TextValue 'D45' '20.43'
Set-TextValue 'D46' '0.997'
Set-TextValue 'D47' '4.93'
Set-TextValue 'E47' '  -2.52%  '
Set-TextValue 'E48' '  -3.21%  '
Set-TextValue 'E49' '  -0.98%  '
Set-TextValue 'D50' '0.0920'
Set-TextValue 'E50' '  -3.51%  '
Set-TextValue 'D51' '18.28'
Set-TextValue 'E51' '  -7.46%  '
